$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update task statuses (D column)
$ws.Range("D5").Value = "In progress"
$ws.Range("D14").Value = "Done"
$ws.Range("D15").Value = "In progress"

# Remove the now-obsolete "TBD" link row (row 31: crackstation hashing link)
$ws.Rows("31:31").Delete()

# Reset sheet view: scroll back to top-left and change selection to E6
$ws.Range("E6").Select() | Out-Null
